$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 254, pushing existing rows 254-289 down to 255-290.
$ws.Rows.Item(254).Insert()

# Populate the newly inserted row 254 with the new data record.
$ws.Range("A254").Value = 10
$ws.Range("B254").Value = "Vega Modelo de Temuco"
$ws.Range("C254").Value = "La Araucanía"
$ws.Range("D254").Value = [DateTime]"2021-11-05"
$ws.Range("E254").Value = 9
$ws.Range("F254").Value = "Fruta"
$ws.Range("G254").Value = 100103
$ws.Range("H254").Value = "Frutos de hueso (carozo)"
$ws.Range("I254").Value = 100103006
$ws.Range("J254").Value = "Nectarín"
$ws.Range("K254").Value = "Early Glo"
$ws.Range("L254").Value = "Primera"
$ws.Range("M254").Value = 55
$ws.Range("N254").Value = 32000
$ws.Range("O254").Value = 32000
$ws.Range("P254").Value = 32000
$ws.Range("Q254").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R254").Value = "Provincia de Limarí"
$ws.Range("S254").Value = 1778
$ws.Range("T254").Value = 18
